# Generate Report for Handback
#
# Records a "handback" pass on the localization status report:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every tracked file.
#   - A new pair of columns is populated/linked: "Latest Target File" (E)
#     and "Latest Handback File" (F), mirroring the existing source (A)
#     and handoff (C) links.
#   - "Latest Handback DateTime" (G) is stamped with the handback time.
# This is applied identically to both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

# Custom hyperlink color already used in this workbook (RGB FF6495ED /
# cornflowerblue), expressed as an OLE (BGR) value for Font.Color.
$hyperlinkColor = 15570276

function Set-HandbackLink($ws, $cellRef, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText) | Out-Null
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

function Update-LanguageSheet($sheetName, $mdUrl, $mdDisplay, $xlfUrl, $xlfDisplay, $md2Url, $md2Display, $xlf2Url, $xlf2Display, $handbackTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # New "Latest Target File" (E) / "Latest Handback File" (F) columns,
    # mirroring the existing Markdown source (A) / xlf handoff (C) links.
    Set-HandbackLink $ws "E2" $mdUrl $mdDisplay
    Set-HandbackLink $ws "F2" $xlfUrl $xlfDisplay
    Set-HandbackLink $ws "E3" $md2Url $md2Display
    Set-HandbackLink $ws "F3" $xlf2Url $xlf2Display

    # Latest Handback DateTime (G2/G3) gets a fresh timestamp.
    $ws.Range("G2").Value = $handbackTime
    $ws.Range("G3").Value = $handbackTime
}

Update-LanguageSheet "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e33b08e25774e95864db62a23b5a148b363a47bd/e2e/aacb56dc-0b34-46ca-b76a-52dcf60615ef.md" `
    "aacb56dc-0b34-46ca-b76a-52dcf60615ef.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/450c007025b5557573ce6dada3d88b29c2a3a823/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/aacb56dc-0b34-46ca-b76a-52dcf60615ef.bc46d157e19ad3a880572edd63052e96973858f2.zh-cn.xlf" `
    "aacb56dc-0b34-46ca-b76a-52dcf60615ef.bc46d157e19ad3a880572edd63052e96973858f2.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e33b08e25774e95864db62a23b5a148b363a47bd/e2e/f4df2f73-d872-47a0-95be-1856bc9aac1d.md" `
    "f4df2f73-d872-47a0-95be-1856bc9aac1d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/450c007025b5557573ce6dada3d88b29c2a3a823/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/f4df2f73-d872-47a0-95be-1856bc9aac1d.29feef96d61737d1a5a4a0597906315098594ffb.zh-cn.xlf" `
    "f4df2f73-d872-47a0-95be-1856bc9aac1d.29feef96d61737d1a5a4a0597906315098594ffb.zh-cn.xlf" `
    "2016-02-26 07:13:57"

Update-LanguageSheet "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e33b08e25774e95864db62a23b5a148b363a47bd/e2e/aacb56dc-0b34-46ca-b76a-52dcf60615ef.md" `
    "aacb56dc-0b34-46ca-b76a-52dcf60615ef.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ebf527794bf7579889a8655d556f769af8b6b7e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/aacb56dc-0b34-46ca-b76a-52dcf60615ef.bc46d157e19ad3a880572edd63052e96973858f2.de-de.xlf" `
    "aacb56dc-0b34-46ca-b76a-52dcf60615ef.bc46d157e19ad3a880572edd63052e96973858f2.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/e33b08e25774e95864db62a23b5a148b363a47bd/e2e/f4df2f73-d872-47a0-95be-1856bc9aac1d.md" `
    "f4df2f73-d872-47a0-95be-1856bc9aac1d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ebf527794bf7579889a8655d556f769af8b6b7e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/f4df2f73-d872-47a0-95be-1856bc9aac1d.29feef96d61737d1a5a4a0597906315098594ffb.de-de.xlf" `
    "f4df2f73-d872-47a0-95be-1856bc9aac1d.29feef96d61737d1a5a4a0597906315098594ffb.de-de.xlf" `
    "2016-02-26 07:14:26"
